$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (price) cells to text format so numeric-looking values
# like "1.001" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.992.11"
$ws.Range("E2").Value = "  -0.11%  "

$ws.Range("D3").Value = "1.909.19"
$ws.Range("E3").Value = "  +0.17%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "0.7921"
$ws.Range("E5").Value = "  +4.58%  "

$ws.Range("D6").Value = "241.68"
$ws.Range("E6").Value = "  +0.17%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  +2.63%  "

$ws.Range("D9").Value = "26.28"
$ws.Range("E9").Value = "  +3.19%  "

$ws.Range("D10").Value = "0.06913"
$ws.Range("E10").Value = "  +0.33%  "

$ws.Range("D11").Value = "0.08001"
$ws.Range("E11").Value = "  -0.04%  "

$ws.Range("D12").Value = "1.905.25"
$ws.Range("E12").Value = "  +0.02%  "

$ws.Range("D13").Value = "0.7412"
$ws.Range("E13").Value = "  -1.68%  "

$ws.Range("D14").Value = "5.189"
$ws.Range("E14").Value = "  -1.05%  "

$ws.Range("E15").Value = "  +1.43%  "

$ws.Range("D16").Value = "29.995.70"
$ws.Range("E16").Value = "  -0.10%  "

$ws.Range("D17").Value = "13.96"
$ws.Range("E17").Value = "  -0.28%  "

$ws.Range("D18").Value = "5.867"
$ws.Range("E18").Value = "  -5.07%  "

$ws.Range("D19").Value = "245.91"
$ws.Range("E19").Value = "  +3.76%  "

$ws.Range("D20").Value = "0.000007744"
$ws.Range("E20").Value = "  +0.36%  "

$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("D22").Value = "2.152.11"
$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").Value = "6.829"
$ws.Range("E24").Value = "  -2.70%  "

$ws.Range("D25").Value = "168.00"
$ws.Range("E25").Value = "  +1.12%  "

$ws.Range("D26").Value = "9.245"
$ws.Range("E26").Value = "  -0.53%  "

$ws.Range("D27").Value = "0.1396"
$ws.Range("E27").Value = "  +7.19%  "

$ws.Range("E28").Value = "  +0.47%  "

$ws.Range("D29").Value = "2.029"
$ws.Range("E29").Value = "  -1.81%  "

$ws.Range("E30").Value = "  +1.57%  "

$ws.Range("D31").Value = "1.512"
$ws.Range("E31").Value = "  -0.71%  "

$ws.Range("D32").Value = "4.309"
$ws.Range("E32").Value = "  +0.10%  "

$ws.Range("D33").Value = "4.085"
$ws.Range("E33").Value = "  +1.02%  "

$ws.Range("D34").Value = "0.05541"

$ws.Range("E35").Value = "  -2.35%  "

$ws.Range("D36").Value = "0.7319"
$ws.Range("E36").Value = "  -0.59%  "

$ws.Range("D37").Value = "2.720"
$ws.Range("E37").Value = "  +0.24%  "

$ws.Range("D38").Value = "0.01923"
$ws.Range("E38").Value = "  -0.86%  "

$ws.Range("D39").Value = "2.784"
$ws.Range("E39").Value = "  +0.76%  "

$ws.Range("D40").Value = "6.115"
$ws.Range("E40").Value = "  -1.97%  "

$ws.Range("E41").Value = "  -0.82%  "

$ws.Range("D42").Value = "72.26"
$ws.Range("E42").Value = "  -0.66%  "

$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").Value = "0.8339"
$ws.Range("E44").Value = "  +0.35%  "

$ws.Range("D45").Value = "1.877"
$ws.Range("E45").Value = "  -3.32%  "

$ws.Range("D46").Value = "100.53"
$ws.Range("E46").Value = "  -0.94%  "

$ws.Range("D47").Value = "7.530"
$ws.Range("E47").Value = "  -1.60%  "

$ws.Range("D48").Value = "987.33"
$ws.Range("E48").Value = "  +7.61%  "

$ws.Range("D49").Value = "2.057.99"
$ws.Range("E49").Value = "  +0.03%  "

$ws.Range("D50").Value = "36.19"
$ws.Range("E50").Value = "  -0.85%  "

$ws.Range("D51").Value = "2.804"
$ws.Range("E51").Value = "  +6.60%  "

# Restore the default (Normal) cell style on column D so no stray
# style/number-format attribute is left behind on the cells.
$ws.Range("D2:D51").Style = "Normal"